# Restored from revision of admin on 01/12/2021 04:19:46 PM.TEST Author: admin. Type: SAVE.
# The only substantive data change in this revision is the "Integer min"
# value for rule R30 (row 10, column C) on the "Rules" sheet: 18 -> 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
